# Autogenerated on Sun Feb 01 2015 22:24:41 GMT-0500 (Eastern Standard Time)
#
# Restructures the Senegal "Data" sheet into a "Summary" sheet: a new
# "Source Type" sub-section is added above the existing Micro/SMEs/MSMEs
# table, and a second sub-section ("Value added to the economy") plus a
# sources/citation block are appended below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename sheet: Data -> Summary -------------------------------------
$ws.Name = "Summary"

# --- clear the old table (it is being replaced / relocated) -----------
$ws.Range("A5:D7").Clear()

# --- new sub-section header above the first table ----------------------
$ws.Range("A9").Value = "Source Type: SME Associations (Most Widely Used)"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Underline = $true

# --- table 1: Micro / SMEs / MSMEs headers ------------------------------
$ws.Range("B11").Value = "Micro"
$ws.Range("B11").Font.Bold = $true
$ws.Range("C11").Value = "SMEs"
$ws.Range("C11").Font.Bold = $true
$ws.Range("D11").Value = "MSMEs"
$ws.Range("D11").Font.Bold = $true

# Employment (% of total) row
$ws.Range("A12").Value = "Employment (% of total)"
$ws.Range("A12").Font.Bold = $true
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "30"

# Enterprises (% of total) row
$ws.Range("A13").Value = "Enterprises (% of total)"
$ws.Range("A13").Font.Bold = $true
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "90"

# source note for table 1
$ws.Range("A14").Value = "Source: Min PME - ADEPME, 2010"
$ws.Range("A14").Font.Italic = $true

# --- table 2: Micro / SMEs / MSMEs headers (repeated) -------------------
$ws.Range("B16").Value = "Micro"
$ws.Range("B16").Font.Bold = $true
$ws.Range("C16").Value = "SMEs"
$ws.Range("C16").Font.Bold = $true
$ws.Range("D16").Value = "MSMEs"
$ws.Range("D16").Font.Bold = $true

# Value added to the economy (% of total) row
$ws.Range("A17").Value = "Value added to the economy (% of total)"
$ws.Range("A17").Font.Bold = $true
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20"

# source note for table 2
$ws.Range("A18").Value = "Source: Min PME - ADEPME, 2010"
$ws.Range("A18").Font.Italic = $true

# --- sources / citation block -------------------------------------------
$ws.Range("A23").Value = "Min PME - ADEPME"
$ws.Range("A23").Font.Bold = $true

$ws.Range("A24").Value = "Ministere des mines, de l'Industrie, de l'Agro-industrie et des PME, Direction des Petites et Moyennes Entreprises (Min PME - ADEPME), ""LETTRE  DE  POLITIQUE SECTORIELLE  DES  PME"", 2010, p. 9. Available at http://www.senegal-entreprises.net/3-download/lettre-politique-sectorielle-10-2010.pdf"
$ws.Range("A24").Font.Italic = $true
